$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("S2").Formula = '=HYPERLINK("https://klasma.github.io/LoggingDetectiveFiles/Logging_MONSTERAS/artfynd/A 32298-2023.xlsx"; "A 32298-2023")'
